$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the old hyperlink definition and the inherited "@ text" column /
# tinted header formatting so the sheet starts from a clean slate.
$ws.Hyperlinks.Delete()
$ws.Range("A:B").ClearFormats()

# --- Values ---------------------------------------------------------------
$ws.Range("A1").Value = "KEY"
$ws.Range("B1").Value = "Value"

$ws.Range("A2").Value = "UK"
$ws.Range("B2").Value = "Birmingham, Manchester"

$ws.Range("A3").Value = "password"
$ws.Range("B3").Value = "Tuan@728"

$ws.Range("A4").Value = "username"
$ws.Range("B4").Value = 906249919

$ws.Range("A5").Value = "mail"

# --- Header style: accent1 themed fill (no tint) ---------------------------
$headerRange = $ws.Range("A1:B1")
$headerRange.Interior.ThemeColor = 5
$headerRange.Interior.TintAndShade = 0

# --- Hyperlink for the password value --------------------------------------
$ws.Hyperlinks.Add($ws.Range("B3"), "mailto:Tuan@728") | Out-Null

# --- Column widths ----------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 8.43
$ws.Columns.Item(2).ColumnWidth = 28.15

# --- Selection / view state --------------------------------------------------
$ws.Range("B5").Select()
